# Petrol Truly Unlimited - Testing.xlsx edit
# Fixes the C4 typo ("fuelling" -> "fueling"), fills in test-plan rows 22-30
# (rows 22-24 already had matching text; rows 25-30 were previously blank),
# adjusts a few row heights, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: fix the duplicate/typo string used for "Purpose of test" ---
$ws.Range("C4").Value = "Close program after at least one vehicle finnishes fueling"

# --- Rows 22-24: re-affirm existing content (values already correct) ---
$ws.Range("B22").Value = "Vehicles dont take more than the maximum fueling time"
$ws.Range("C22").Value = "See if any vehicle takes more than 18 seconds to fuel (default max time)"
$ws.Range("D22").Value = "Vehicles take 18 seconds or less to fuel"
$ws.Range("E22").Value = "Pass"

$ws.Range("B23").Value = "Vehicles can leave after fueling"
$ws.Range("C23").Value = "After a vehicle fuels he can leave"
$ws.Range("D23").Value = "The vehicle leavs the pump making it available"
$ws.Range("E23").Value = "Pass"

$ws.Range("B24").Value = "The pumps are dispensing the correct type of fuel"
$ws.Range("C24").Value = "When a vehicle enters the correct fuel is dispensed"
$ws.Range("D24").Value = "Pumps dispense the correct fuel"
$ws.Range("E24").Value = "Pass"

# --- Row 25: new test case ---
$ws.Range("B25").Value = "After the vehicle leaves the receipt is saved"
$ws.Range("C25").Value = "Wait for a vehicle leaves the pump and look for a receipt"
$ws.Range("D25").Value = "A new receipt is saved"
$ws.Range("E25").Value = "Pass"
$ws.Rows.Item(25).RowHeight = 15

# --- Row 26: new test case ---
$ws.Range("B26").Value = "Pumps information is updated"
$ws.Range("C26").Value = "Move the mouse over the pump to see the statistics of each pump"
$ws.Range("D26").Value = "Everytime a vehicle enters or leaved the stats is updated"
$ws.Range("E26").Value = "Pass"
$ws.Rows.Item(26).RowHeight = 25.5

# --- Row 27: new test case ---
$ws.Range("B27").Value = "Each car is spawned with random information"
$ws.Range("C27").Value = "Move the mouse over the vehicle to see the information of each vehicle"
$ws.Range("D27").Value = "A popup should appear with the vehicle information"
$ws.Range("E27").Value = "Pass"

# --- Row 28: new test case ---
$ws.Range("B28").Value = "The queue should have information updated everytime"
$ws.Range("C28").Value = "Move the mouse over the queue to see the statistics"
$ws.Range("D28").Value = "A popup should appear with the queue information"
$ws.Range("E28").Value = "Pass"
$ws.Rows.Item(28).AutoFit()

# --- Row 29: new test case ---
$ws.Range("B29").Value = "Each car in the pump should open a popup with information about the vehicle"
$ws.Range("C29").Value = "Move the mouse over the car that is in the pump"
$ws.Range("D29").Value = "A popup should appear with the vehicle information"
$ws.Range("E29").Value = "Pass"
$ws.Rows.Item(29).RowHeight = 25.5

# --- Row 30: new test case ---
$ws.Range("B30").Value = "Should be possible to scroll through the last 200 receipts"
$ws.Range("C30").Value = "When the receipt list has enough items a scroll bar should appear and we it should be able to scroll"
$ws.Range("D30").Value = "List should be able to scroll"
$ws.Range("E30").Value = "Pass"

# --- Selection / view state ---
$ws.Activate()
$ws.Range("E30").Select()
